$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" -----------------
# "Overview" sheet: status values live in columns E (zh-cn) and F (de-de),
# rows 2-3.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# Per-locale sheets: status values live in column C, rows 2-3.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Column width change ----------------------------------------------------
# The status columns got narrower (the new text is shorter than the old),
# matching an autofit-style resize from ~17.22 down to ~13.41 "characters".
# Apply the narrower width to the same columns on every sheet. Columns are
# referenced by numeric index (not letter) for reliability.
$newWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth

Write-Host "Done applying localization-status report updates"
